$wb = $excel.ActiveWorkbook

# Sheet 1: GNSP Template -- selection change only
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A2").Select()

# Sheet 2: met_metadata -- zoom + selection change
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$excel.ActiveWindow.Zoom = 173
$ws2.Range("D2").Select()

# Sheet 3: batch_example -- clear column styles, update data, selection change (last = active sheet)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Activate()
$ws3.Columns.EntireColumn.ClearFormats()
$ws3.Range("D2").Value = "LC-ESI"
$ws3.Range("D3").Value = "LC-ESI"
$ws3.Range("D4").Select()
